$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.296.65"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").Value = "'1.804.85"
$ws.Range("E3").Value = "  +2.98%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'339.83"
$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").Value = "'0.4651"
$ws.Range("E7").Value = "  +20.13%  "

$ws.Range("D8").Value = "'0.3814"
$ws.Range("E8").Value = "  +12.65%  "

$ws.Range("D9").Value = "'45.37"
$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("D10").Value = "'1.159"
$ws.Range("E10").Value = "  +4.06%  "

$ws.Range("D11").Value = "'0.07598"
$ws.Range("E11").Value = "  +5.36%  "

$ws.Range("D12").Value = "'22.48"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").Value = "'1.002"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").Value = "'6.349"
$ws.Range("E14").Value = "  +2.66%  "

$ws.Range("D15").Value = "'7.541"
$ws.Range("E15").Value = "  +6.20%  "

$ws.Range("D16").Value = "'1.808.79"
$ws.Range("E16").Value = "  +3.37%  "

$ws.Range("D17").Value = "'0.00001096"
$ws.Range("E17").Value = "  +3.47%  "

$ws.Range("D18").Value = "'0.06735"
$ws.Range("E18").Value = "  +2.04%  "

$ws.Range("D19").Value = "'81.65"
$ws.Range("E19").Value = "  +2.73%  "

$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("E21").Value = "  +4.27%  "

$ws.Range("D22").Value = "'6.437"
$ws.Range("E22").Value = "  +4.08%  "

$ws.Range("D23").Value = "'28.278.05"
$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("D24").Value = "'11.92"
$ws.Range("E24").Value = "  +2.13%  "

$ws.Range("D25").Value = "'2.432"

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'20.68"
$ws.Range("E26").Value = "  +4.08%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'153.65"
$ws.Range("E27").Value = "  -0.69%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.377"
$ws.Range("E28").Value = "  +3.21%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "'2.013.83"
$ws.Range("E29").Value = "  +3.01%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'133.03"
$ws.Range("E30").Value = "  +1.71%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.253"
$ws.Range("E31").Value = "  -3.20%  "

$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "'4.041"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.09598"
$ws.Range("E33").Value = "  +9.19%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.852"
$ws.Range("E34").Value = "  +0.58%  "

$ws.Range("B35").Value = "Algorand"
$ws.Range("C35").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D35").Value = "'0.2314"
$ws.Range("E35").Value = "  +9.79%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'12.10"
$ws.Range("E36").Value = "  -0.55%  "

$ws.Range("E37").Value = "  +3.57%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06357"
$ws.Range("E38").Value = "  +3.97%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.278"
$ws.Range("E39").Value = "  +2.70%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6629"
$ws.Range("E40").Value = "  +1.35%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.243"
$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'8.397"
$ws.Range("E42").Value = "  +4.66%  "

$ws.Range("E43").Value = "  -3.07%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.21"
$ws.Range("E44").Value = "  +4.15%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6143"
$ws.Range("E46").Value = "  +1.66%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.872"
$ws.Range("E47").Value = "  +1.40%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'130.94"
$ws.Range("E48").Value = "  +3.16%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.041"
$ws.Range("E49").Value = "  +2.36%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07168"
$ws.Range("E50").Value = "  +2.93%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.179"
$ws.Range("E51").Value = "  +1.54%  "
